$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.948.59"
$ws.Range("E2").Value = "  -3.97%  "

$ws.Range("D3").Value = "3.490.99"
$ws.Range("E3").Value = "  -5.64%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").Value = "3.481.80"
$ws.Range("E8").Value = "  -5.60%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.187"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.65"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.591"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000274"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "671.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "4.054.53"
$ws.Range("E16").Value = "  -5.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.47%  "

$ws.Range("D18").Value = "68.921.39"
$ws.Range("E18").Value = "  -4.14%  "

$ws.Range("D19").Value = "3.497.01"
$ws.Range("E19").Value = "  -5.51%  "

$ws.Range("E20").Value = "  -1.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.895"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -10.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.18"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.70"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -8.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -9.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.20"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "595.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.06%  "

$ws.Range("B36").Value = "Cosmos"
$ws.Range("C36").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.80"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.41%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -15.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.87"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0434"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.332"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.78%  "

$ws.Range("D43").Value = "3.394.06"
$ws.Range("E43").Value = "  -9.79%  "

$ws.Range("E44").Value = "  -6.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.03"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.66%  "

$ws.Range("D46").Value = "0.0₃0700"
$ws.Range("E46").Value = "  -9.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.87"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.57"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.76%  "

$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.73"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +16.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.11%  "
